# Update "想去人数" (interested-count) figures to the newly scraped values.
# Affects both the per-category sheet ("展览") and the combined "全部类型"
# sheet, which both contain the same three events.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - rows 2,3,4 hold the three events in question.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 409
$wsExhibit.Range("F3").Value = 2347
$wsExhibit.Range("F4").Value = 113

# Sheet "全部类型" (All types) - combined listing; same events appear as
# rows 2, 7, 8.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 409
$wsAll.Range("F7").Value = 2347
$wsAll.Range("F8").Value = 113
